$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model")
$ws.Rows(2).Delete()
